# November.xlsx — "data update 12.9 with previous commitment push failure re-push"
# Append five new ledger rows (237-242, i.e. entries #236-#241 plus a subtotal
# row) to Sheet1, matching the style/format of the nearby existing rows, and
# move the active selection from E232 to E240 (the new subtotal row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Carry formatting forward from the closest existing rows of the same
#     kind, so new cells inherit the right number formats / fills instead of
#     plain defaults. xlPasteFormats = -4122
$ws.Range("A236:F236").Copy()
$ws.Range("A237:F239").PasteSpecial(-4122)

$ws.Range("A232:H232").Copy()
$ws.Range("A240:H240").PasteSpecial(-4122)

$ws.Range("A233:F234").Copy()
$ws.Range("A241:F242").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 237 (entry #236): 打车 / 从学校到澡堂
$ws.Cells.Item(237, 2).Value = 45634
$ws.Cells.Item(237, 3).Value = "打车"
$ws.Cells.Item(237, 4).Value = -7.49
$ws.Cells.Item(237, 5).Value = "从学校到澡堂"
$ws.Cells.Item(237, 6).Formula = "=F236+D237"

# --- Row 238 (entry #237): M / 吃吃M记
$ws.Cells.Item(238, 2).Value = 45634
$ws.Cells.Item(238, 3).Value = "M"
$ws.Cells.Item(238, 4).Value = -27.8
$ws.Cells.Item(238, 5).Value = "吃吃M记"
$ws.Cells.Item(238, 6).Formula = "=F237+D238"

# --- Row 239 (entry #238): 打车 / 从澡堂旁的M回学校
$ws.Cells.Item(239, 2).Value = 45634
$ws.Cells.Item(239, 3).Value = "打车"
$ws.Cells.Item(239, 4).Value = -7.47
$ws.Cells.Item(239, 5).Value = "从澡堂旁的M回学校"
$ws.Cells.Item(239, 6).Formula = "=F238+D239"

# --- Row 240 (entry #239): 小结 subtotal for the 45634 day (rows 236-239)
$ws.Cells.Item(240, 2).Value = 45634
$ws.Cells.Item(240, 3).Value = "小结"
$ws.Cells.Item(240, 4).Formula = "=SUM(D236:D239)"
$ws.Cells.Item(240, 5).Value = "*"
$ws.Cells.Item(240, 6).Formula = "=F239"

# --- Row 241 (entry #240): 早饭 / 经典套餐
$ws.Cells.Item(241, 2).Value = 45635
$ws.Cells.Item(241, 3).Value = "早饭"
$ws.Cells.Item(241, 4).Value = -9.1
$ws.Cells.Item(241, 5).Value = "经典套餐"
$ws.Cells.Item(241, 6).Formula = "=F240+D241"

# --- Row 242 (entry #241): 迅雷会员 / 两个月送航旅月卡
$ws.Cells.Item(242, 2).Value = 45635
$ws.Cells.Item(242, 3).Value = "迅雷会员"
$ws.Cells.Item(242, 4).Value = -30
$ws.Cells.Item(242, 5).Value = "两个月送航旅月卡"
$ws.Cells.Item(242, 6).Formula = "=F241+D242"

# --- Move the active selection to the new subtotal cell, as in the saved file
$ws.Range("E240").Select() | Out-Null
